$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row swaps / rotations (update F:V in place) ---
# Row 20
$ws.Cells.Item(20, 6).Value = 'Hellerup'
$ws.Cells.Item(20, 7).Value = 2
$ws.Cells.Item(20, 8).Value = 'Nykobing'
$ws.Cells.Item(20, 9).Value = 3
$ws.Cells.Item(20, 10).Value = 2.33
$ws.Cells.Item(20, 11).Value = '25/08/2023 11:42'
$ws.Cells.Item(20, 12).Value = 2.41
$ws.Cells.Item(20, 13).Value = '25/08/2023 18:44'
$ws.Cells.Item(20, 14).Value = 3.42
$ws.Cells.Item(20, 15).Value = '25/08/2023 11:42'
$ws.Cells.Item(20, 16).Value = 3.73
$ws.Cells.Item(20, 17).Value = '25/08/2023 17:25'
$ws.Cells.Item(20, 18).Value = 2.81
$ws.Cells.Item(20, 19).Value = '25/08/2023 11:42'
$ws.Cells.Item(20, 20).Value = 2.57
$ws.Cells.Item(20, 21).Value = '25/08/2023 18:44'
$ws.Cells.Item(20, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/hellerup-nykobing/EqhUX5vt/'

# Row 21
$ws.Cells.Item(21, 6).Value = 'FA 2000'
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 'AB Copenhagen'
$ws.Cells.Item(21, 9).Value = 2
$ws.Cells.Item(21, 10).Value = 3.11
$ws.Cells.Item(21, 11).Value = '25/08/2023 11:42'
$ws.Cells.Item(21, 12).Value = 3.13
$ws.Cells.Item(21, 13).Value = '25/08/2023 18:44'
$ws.Cells.Item(21, 14).Value = 3.56
$ws.Cells.Item(21, 15).Value = '25/08/2023 11:42'
$ws.Cells.Item(21, 16).Value = 3.75
$ws.Cells.Item(21, 17).Value = '25/08/2023 18:44'
$ws.Cells.Item(21, 18).Value = 2.06
$ws.Cells.Item(21, 19).Value = '25/08/2023 11:42'
$ws.Cells.Item(21, 20).Value = 2.06
$ws.Cells.Item(21, 21).Value = '25/08/2023 18:44'
$ws.Cells.Item(21, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/frederiksberg-alliancen-2000-ab-copenhagen/lMQNCM8U/'

# Row 35
$ws.Cells.Item(35, 6).Value = 'Esbjerg'
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = 'Nykobing'
$ws.Cells.Item(35, 9).Value = 4
$ws.Cells.Item(35, 10).Value = 1.36
$ws.Cells.Item(35, 11).Value = '10/09/2023 09:12'
$ws.Cells.Item(35, 12).Value = 1.32
$ws.Cells.Item(35, 13).Value = '10/09/2023 13:24'
$ws.Cells.Item(35, 14).Value = 5.13
$ws.Cells.Item(35, 15).Value = '10/09/2023 09:12'
$ws.Cells.Item(35, 16).Value = 5.56
$ws.Cells.Item(35, 17).Value = '10/09/2023 13:24'
$ws.Cells.Item(35, 18).Value = 6.23
$ws.Cells.Item(35, 19).Value = '10/09/2023 09:12'
$ws.Cells.Item(35, 20).Value = 7.4
$ws.Cells.Item(35, 21).Value = '10/09/2023 13:24'
$ws.Cells.Item(35, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/esbjerg-nykobing/EwkKt1NA/'

# Row 36
$ws.Cells.Item(36, 6).Value = 'Aarhus Fremad'
$ws.Cells.Item(36, 7).Value = 3
$ws.Cells.Item(36, 8).Value = 'F. Amager'
$ws.Cells.Item(36, 9).Value = 1
$ws.Cells.Item(36, 10).Value = 1.34
$ws.Cells.Item(36, 11).Value = '09/09/2023 02:12'
$ws.Cells.Item(36, 12).Value = 1.49
$ws.Cells.Item(36, 13).Value = '10/09/2023 13:31'
$ws.Cells.Item(36, 14).Value = 4.74
$ws.Cells.Item(36, 15).Value = '09/09/2023 02:12'
$ws.Cells.Item(36, 16).Value = 4.45
$ws.Cells.Item(36, 17).Value = '10/09/2023 13:31'
$ws.Cells.Item(36, 18).Value = 6.73
$ws.Cells.Item(36, 19).Value = '09/09/2023 02:12'
$ws.Cells.Item(36, 20).Value = 5.68
$ws.Cells.Item(36, 21).Value = '10/09/2023 13:31'
$ws.Cells.Item(36, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-fremad-amager/hUmOuLxH/'

# Row 55
$ws.Cells.Item(55, 6).Value = 'Aarhus Fremad'
$ws.Cells.Item(55, 7).Value = 1
$ws.Cells.Item(55, 8).Value = 'AB Copenhagen'
$ws.Cells.Item(55, 9).Value = 1
$ws.Cells.Item(55, 10).Value = 1.55
$ws.Cells.Item(55, 11).Value = '06/10/2023 01:13'
$ws.Cells.Item(55, 12).Value = 1.48
$ws.Cells.Item(55, 13).Value = '07/10/2023 13:58'
$ws.Cells.Item(55, 14).Value = 4.09
$ws.Cells.Item(55, 15).Value = '06/10/2023 01:13'
$ws.Cells.Item(55, 16).Value = 4.81
$ws.Cells.Item(55, 17).Value = '07/10/2023 13:58'
$ws.Cells.Item(55, 18).Value = 4.41
$ws.Cells.Item(55, 19).Value = '06/10/2023 01:13'
$ws.Cells.Item(55, 20).Value = 5.36
$ws.Cells.Item(55, 21).Value = '07/10/2023 13:58'
$ws.Cells.Item(55, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-ab-copenhagen/AXvIm9ze/'

# Row 56
$ws.Cells.Item(56, 6).Value = 'Middelfart'
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 'Roskilde'
$ws.Cells.Item(56, 9).Value = 1
$ws.Cells.Item(56, 10).Value = 2.18
$ws.Cells.Item(56, 11).Value = '06/10/2023 01:12'
$ws.Cells.Item(56, 12).Value = 2.26
$ws.Cells.Item(56, 13).Value = '07/10/2023 10:07'
$ws.Cells.Item(56, 14).Value = 3.21
$ws.Cells.Item(56, 15).Value = '06/10/2023 01:12'
$ws.Cells.Item(56, 16).Value = 3.32
$ws.Cells.Item(56, 17).Value = '07/10/2023 12:01'
$ws.Cells.Item(56, 18).Value = 2.88
$ws.Cells.Item(56, 19).Value = '06/10/2023 01:12'
$ws.Cells.Item(56, 20).Value = 3.01
$ws.Cells.Item(56, 21).Value = '07/10/2023 10:07'
$ws.Cells.Item(56, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/middelfart-roskilde/S8kDlkLl/'

# Row 57
$ws.Cells.Item(57, 6).Value = 'Skive'
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 'Nykobing'
$ws.Cells.Item(57, 9).Value = 2
$ws.Cells.Item(57, 10).Value = 2.58
$ws.Cells.Item(57, 11).Value = '06/10/2023 01:13'
$ws.Cells.Item(57, 12).Value = 3.55
$ws.Cells.Item(57, 13).Value = '07/10/2023 13:40'
$ws.Cells.Item(57, 14).Value = 3.24
$ws.Cells.Item(57, 15).Value = '06/10/2023 01:13'
$ws.Cells.Item(57, 16).Value = 3.55
$ws.Cells.Item(57, 17).Value = '07/10/2023 13:40'
$ws.Cells.Item(57, 18).Value = 2.39
$ws.Cells.Item(57, 19).Value = '06/10/2023 01:13'
$ws.Cells.Item(57, 20).Value = 1.97
$ws.Cells.Item(57, 21).Value = '07/10/2023 13:40'
$ws.Cells.Item(57, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/skive-nykobing/lIj9kV5r/'

# Row 67
$ws.Cells.Item(67, 6).Value = 'Middelfart'
$ws.Cells.Item(67, 7).Value = 2
$ws.Cells.Item(67, 8).Value = 'Nykobing'
$ws.Cells.Item(67, 9).Value = 2
$ws.Cells.Item(67, 10).Value = 1.95
$ws.Cells.Item(67, 11).Value = '20/10/2023 01:13'
$ws.Cells.Item(67, 12).Value = 2.1
$ws.Cells.Item(67, 13).Value = '21/10/2023 13:41'
$ws.Cells.Item(67, 14).Value = 3.43
$ws.Cells.Item(67, 15).Value = '20/10/2023 01:13'
$ws.Cells.Item(67, 16).Value = 3.66
$ws.Cells.Item(67, 17).Value = '21/10/2023 13:41'
$ws.Cells.Item(67, 18).Value = 3.21
$ws.Cells.Item(67, 19).Value = '20/10/2023 01:13'
$ws.Cells.Item(67, 20).Value = 3.12
$ws.Cells.Item(67, 21).Value = '21/10/2023 13:41'
$ws.Cells.Item(67, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/middelfart-nykobing/CzzEZORI/'

# Row 68
$ws.Cells.Item(68, 6).Value = 'Brabrand'
$ws.Cells.Item(68, 7).Value = 2
$ws.Cells.Item(68, 8).Value = 'AB Copenhagen'
$ws.Cells.Item(68, 9).Value = 2
$ws.Cells.Item(68, 10).Value = 3.31
$ws.Cells.Item(68, 11).Value = '20/10/2023 01:13'
$ws.Cells.Item(68, 12).Value = 3.55
$ws.Cells.Item(68, 13).Value = '21/10/2023 13:41'
$ws.Cells.Item(68, 14).Value = 3.44
$ws.Cells.Item(68, 15).Value = '20/10/2023 01:13'
$ws.Cells.Item(68, 16).Value = 3.49
$ws.Cells.Item(68, 17).Value = '21/10/2023 13:41'
$ws.Cells.Item(68, 18).Value = 1.95
$ws.Cells.Item(68, 19).Value = '20/10/2023 01:13'
$ws.Cells.Item(68, 20).Value = 1.99
$ws.Cells.Item(68, 21).Value = '21/10/2023 13:41'
$ws.Cells.Item(68, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/brabrand-ab-copenhagen/Iys5yoc6/'

# Row 69
$ws.Cells.Item(69, 6).Value = 'Esbjerg'
$ws.Cells.Item(69, 7).Value = 3
$ws.Cells.Item(69, 8).Value = 'FA 2000'
$ws.Cells.Item(69, 9).Value = 1
$ws.Cells.Item(69, 10).Value = 1.21
$ws.Cells.Item(69, 11).Value = '20/10/2023 01:13'
$ws.Cells.Item(69, 12).Value = 1.18
$ws.Cells.Item(69, 13).Value = '21/10/2023 06:38'
$ws.Cells.Item(69, 14).Value = 6.19
$ws.Cells.Item(69, 15).Value = '20/10/2023 01:13'
$ws.Cells.Item(69, 16).Value = 7.26
$ws.Cells.Item(69, 17).Value = '21/10/2023 13:25'
$ws.Cells.Item(69, 18).Value = 7.96
$ws.Cells.Item(69, 19).Value = '20/10/2023 01:13'
$ws.Cells.Item(69, 20).Value = 11.24
$ws.Cells.Item(69, 21).Value = '21/10/2023 13:25'
$ws.Cells.Item(69, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/esbjerg-frederiksberg-alliancen-2000/00W9z5CC/'

# --- New rows 73-78 ---
# Row 73
$ws.Cells.Item(72, 1).Copy($ws.Cells.Item(73, 1))
$ws.Cells.Item(72, 5).Copy($ws.Cells.Item(73, 5))
$ws.Cells.Item(73, 1).Value = 72
$ws.Cells.Item(73, 2).Value = 'denmark'
$ws.Cells.Item(73, 3).Value = '2nd-division'
$ws.Cells.Item(73, 4).Value = '2023-2024'
$ws.Cells.Item(73, 5).Value = 45226.79166666666
$ws.Cells.Item(73, 6).Value = 'Hellerup'
$ws.Cells.Item(73, 7).Value = 3
$ws.Cells.Item(73, 8).Value = 'F. Amager'
$ws.Cells.Item(73, 9).Value = 1
$ws.Cells.Item(73, 10).Value = 2.59
$ws.Cells.Item(73, 11).Value = '26/10/2023 07:12'
$ws.Cells.Item(73, 12).Value = 3.2
$ws.Cells.Item(73, 13).Value = '27/10/2023 18:23'
$ws.Cells.Item(73, 14).Value = 3.41
$ws.Cells.Item(73, 15).Value = '26/10/2023 07:12'
$ws.Cells.Item(73, 16).Value = 3.52
$ws.Cells.Item(73, 17).Value = '27/10/2023 18:57'
$ws.Cells.Item(73, 18).Value = 2.3
$ws.Cells.Item(73, 19).Value = '26/10/2023 07:12'
$ws.Cells.Item(73, 20).Value = 2.11
$ws.Cells.Item(73, 21).Value = '27/10/2023 18:57'
$ws.Cells.Item(73, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/hellerup-fremad-amager/pO4ShnSt/'

# Row 74
$ws.Cells.Item(72, 1).Copy($ws.Cells.Item(74, 1))
$ws.Cells.Item(72, 5).Copy($ws.Cells.Item(74, 5))
$ws.Cells.Item(74, 1).Value = 73
$ws.Cells.Item(74, 2).Value = 'denmark'
$ws.Cells.Item(74, 3).Value = '2nd-division'
$ws.Cells.Item(74, 4).Value = '2023-2024'
$ws.Cells.Item(74, 5).Value = 45227.57291666666
$ws.Cells.Item(74, 6).Value = 'FA 2000'
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 'Skive'
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 2.18
$ws.Cells.Item(74, 11).Value = '27/10/2023 02:12'
$ws.Cells.Item(74, 12).Value = 2.22
$ws.Cells.Item(74, 13).Value = '28/10/2023 13:26'
$ws.Cells.Item(74, 14).Value = 3.25
$ws.Cells.Item(74, 15).Value = '27/10/2023 02:12'
$ws.Cells.Item(74, 16).Value = 3.39
$ws.Cells.Item(74, 17).Value = '28/10/2023 13:26'
$ws.Cells.Item(74, 18).Value = 2.96
$ws.Cells.Item(74, 19).Value = '27/10/2023 02:12'
$ws.Cells.Item(74, 20).Value = 3.06
$ws.Cells.Item(74, 21).Value = '28/10/2023 13:26'
$ws.Cells.Item(74, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/frederiksberg-alliancen-2000-skive/OE3Wi6sm/'

# Row 75
$ws.Cells.Item(72, 1).Copy($ws.Cells.Item(75, 1))
$ws.Cells.Item(72, 5).Copy($ws.Cells.Item(75, 5))
$ws.Cells.Item(75, 1).Value = 74
$ws.Cells.Item(75, 2).Value = 'denmark'
$ws.Cells.Item(75, 3).Value = '2nd-division'
$ws.Cells.Item(75, 4).Value = '2023-2024'
$ws.Cells.Item(75, 5).Value = 45227.58333333334
$ws.Cells.Item(75, 6).Value = 'Roskilde'
$ws.Cells.Item(75, 7).Value = 4
$ws.Cells.Item(75, 8).Value = 'Brabrand'
$ws.Cells.Item(75, 9).Value = 4
$ws.Cells.Item(75, 10).Value = 1.56
$ws.Cells.Item(75, 11).Value = '27/10/2023 02:12'
$ws.Cells.Item(75, 12).Value = 1.48
$ws.Cells.Item(75, 13).Value = '27/10/2023 10:36'
$ws.Cells.Item(75, 14).Value = 3.92
$ws.Cells.Item(75, 15).Value = '27/10/2023 02:12'
$ws.Cells.Item(75, 16).Value = 4.27
$ws.Cells.Item(75, 17).Value = '28/10/2023 12:03'
$ws.Cells.Item(75, 18).Value = 4.56
$ws.Cells.Item(75, 19).Value = '27/10/2023 02:12'
$ws.Cells.Item(75, 20).Value = 6.09
$ws.Cells.Item(75, 21).Value = '27/10/2023 10:36'
$ws.Cells.Item(75, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/roskilde-brabrand/Wv7vjpCa/'

# Row 76
$ws.Cells.Item(72, 1).Copy($ws.Cells.Item(76, 1))
$ws.Cells.Item(72, 5).Copy($ws.Cells.Item(76, 5))
$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = 'denmark'
$ws.Cells.Item(76, 3).Value = '2nd-division'
$ws.Cells.Item(76, 4).Value = '2023-2024'
$ws.Cells.Item(76, 5).Value = 45227.58333333334
$ws.Cells.Item(76, 6).Value = 'Aarhus Fremad'
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 'Middelfart'
$ws.Cells.Item(76, 9).Value = 2
$ws.Cells.Item(76, 10).Value = 1.73
$ws.Cells.Item(76, 11).Value = '27/10/2023 02:12'
$ws.Cells.Item(76, 12).Value = 1.82
$ws.Cells.Item(76, 13).Value = '28/10/2023 13:51'
$ws.Cells.Item(76, 14).Value = 3.76
$ws.Cells.Item(76, 15).Value = '27/10/2023 02:12'
$ws.Cells.Item(76, 16).Value = 3.73
$ws.Cells.Item(76, 17).Value = '28/10/2023 13:51'
$ws.Cells.Item(76, 18).Value = 3.69
$ws.Cells.Item(76, 19).Value = '27/10/2023 02:12'
$ws.Cells.Item(76, 20).Value = 3.94
$ws.Cells.Item(76, 21).Value = '28/10/2023 13:51'
$ws.Cells.Item(76, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/aarhus-fremad-middelfart/Iy3ziQdg/'

# Row 77
$ws.Cells.Item(72, 1).Copy($ws.Cells.Item(77, 1))
$ws.Cells.Item(72, 5).Copy($ws.Cells.Item(77, 5))
$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = 'denmark'
$ws.Cells.Item(77, 3).Value = '2nd-division'
$ws.Cells.Item(77, 4).Value = '2023-2024'
$ws.Cells.Item(77, 5).Value = 45227.625
$ws.Cells.Item(77, 6).Value = 'Nykobing'
$ws.Cells.Item(77, 7).Value = 0
$ws.Cells.Item(77, 8).Value = 'Esbjerg'
$ws.Cells.Item(77, 9).Value = 3
$ws.Cells.Item(77, 10).Value = 3.87
$ws.Cells.Item(77, 11).Value = '27/10/2023 03:12'
$ws.Cells.Item(77, 12).Value = 4.51
$ws.Cells.Item(77, 13).Value = '28/10/2023 14:34'
$ws.Cells.Item(77, 14).Value = 4.01
$ws.Cells.Item(77, 15).Value = '27/10/2023 03:12'
$ws.Cells.Item(77, 16).Value = 4.37
$ws.Cells.Item(77, 17).Value = '28/10/2023 14:34'
$ws.Cells.Item(77, 18).Value = 1.65
$ws.Cells.Item(77, 19).Value = '27/10/2023 03:12'
$ws.Cells.Item(77, 20).Value = 1.61
$ws.Cells.Item(77, 21).Value = '28/10/2023 09:39'
$ws.Cells.Item(77, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/nykobing-esbjerg/E3dqk4R5/'

# Row 78
$ws.Cells.Item(72, 1).Copy($ws.Cells.Item(78, 1))
$ws.Cells.Item(72, 5).Copy($ws.Cells.Item(78, 5))
$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = 'denmark'
$ws.Cells.Item(78, 3).Value = '2nd-division'
$ws.Cells.Item(78, 4).Value = '2023-2024'
$ws.Cells.Item(78, 5).Value = 45227.625
$ws.Cells.Item(78, 6).Value = 'Thisted FC'
$ws.Cells.Item(78, 7).Value = 0
$ws.Cells.Item(78, 8).Value = 'AB Copenhagen'
$ws.Cells.Item(78, 9).Value = 1
$ws.Cells.Item(78, 10).Value = 2.52
$ws.Cells.Item(78, 11).Value = '27/10/2023 03:12'
$ws.Cells.Item(78, 12).Value = 2.32
$ws.Cells.Item(78, 13).Value = '28/10/2023 14:59'
$ws.Cells.Item(78, 14).Value = 3.33
$ws.Cells.Item(78, 15).Value = '27/10/2023 03:12'
$ws.Cells.Item(78, 16).Value = 3.54
$ws.Cells.Item(78, 17).Value = '28/10/2023 14:59'
$ws.Cells.Item(78, 18).Value = 2.4
$ws.Cells.Item(78, 19).Value = '27/10/2023 03:12'
$ws.Cells.Item(78, 20).Value = 2.79
$ws.Cells.Item(78, 21).Value = '28/10/2023 14:59'
$ws.Cells.Item(78, 22).Value = 'https://www.betexplorer.com/football/denmark/2nd-division/thisted-fc-ab-copenhagen/hfemlOtC/'
